# "Base do 'adicionar novo dispositivo'" - append new device rows to the
# devices table on the active sheet (A: nome, B: tipo/categoria,
# C: valor, D: estado ligado/desligado).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Tv da cozinha
$ws.Range("A5").Value = "Tv da cozinha"
$ws.Range("B5").Value = "Canal 3"
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = $false

# Row 6: Fechadura dos fundos (sem valor numerico, apenas estado)
$ws.Range("A6").Value = "Fechadura dos fundos"
$ws.Range("B6").Value = "fechadura"
$ws.Range("C6").Value = $false

# Row 7: Ar do banheiro
$ws.Range("A7").Value = "Ar do banheiro"
$ws.Range("B7").Value = "climatizadores"
$ws.Range("C7").Value = 23
$ws.Range("D7").Value = $false
